# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "ADRIANA CRISTINA SARMIENTO BLANCO" worker had two overdue-period
# records (rows 17 and 18) that need to be removed from this account
# statement, leaving only the "JOSE GREGORIO RAMIREZ BURGOS" record.
# Deleting the rows (instead of just clearing them) shifts the trailing
# signature block up and lets Excel garbage-collect the now-unused
# shared strings / styles it had used.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two table rows belonging to the worker being dropped from
# the statement (CC 1001971353 - ADRIANA CRISTINA SARMIENTO BLANCO).
$ws.Range("B17:J18").EntireRow.Delete()

# Recompute the summary figures now that only one worker / one period
# remains: total overdue value, worker count, and period count.
$ws.Range("E11").Value2 = 30208
$ws.Range("C13").Value2 = 1
$ws.Range("F13").Value2 = 1

# Column D ("Nombre Trabajador") was sized to fit the longest name in
# the table; with ADRIANA's (longer) name gone it can shrink back down
# to fit the remaining content.
$ws.Columns.Item(4).AutoFit()
